$d = $word.ActiveDocument

function Replace-DocText($old, $new) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND (replace): [$old]"
        return
    }
    $endIdx = $idx + $old.Length
    $rng = $d.Range($idx, $endIdx)
    $rng.Text = $new
}

function Delete-DocText($old) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND (delete): [$old]"
        return
    }
    $endIdx = $idx + $old.Length
    $rng = $d.Range($idx, $endIdx)
    $rng.Delete()
}

# ---- Title ----
Replace-DocText "Quantum Entanglement: Unveiling the Interconnected Universe" "Unraveling the complexities of medicine"

# ---- Author name: "Isaac Newton" -> "Dr" + "." + " Amy Walters" ----
Replace-DocText "Isaac Newton" "Dr. Amy Walters"

# ---- Email: "isaac" -> "awalters@hschool"; keep ".";
#      "newton@scientificdiscovery" + "." + "com" -> "edu" (drop the extra "." and "com") ----
Replace-DocText "isaac" "awalters@hschool"
Replace-DocText "newton@scientificdiscovery" "edu"
Delete-DocText ".com"

# ---- Body paragraph 1 ----
Replace-DocText "In the realm of quantum physics, a mysterious phenomenon known as quantum entanglement defies conventional intuition and challenges our understanding of reality" "Medicine, the field dedicated to preserving life and alleviating suffering, embarks upon a mission as both multifaceted and profound as the human body itself"

Replace-DocText " This enigmatic concept, first proposed by Albert Einstein, reveals the profound interconnectedness of particles, even when separated by vast distances" " Its practitioners, known as physicians, devote themselves to understanding the intricacies of anatomy, scrutinizing the intricate pathways of disease, and extending a healing hand to those in need"

Replace-DocText " As scientists delved deeper into the intricacies of entanglement, they discovered that the properties of these particles become inextricably linked, regardless of the physical separation between them. This remarkable phenomenon has ignited a revolution in our comprehension of the universe, revealing a tapestry of interconnectedness that transcends the boundaries of space and time" " Through its extensive tapestry of specialization and ongoing scientific advancements, medicine weaves together science, compassion, and dedication, while continuously pushing the boundaries of human understanding"

Replace-DocText "Unveiling the enigmatic nature of quantum entanglement has compelled scientists to re-examine fundamental concepts such as locality and causality" "Within this intricate landscape lies a system of interconnected disciplines that seeks to unravel the complexities of the human condition"

Replace-DocText " Locality dictates that no physical influence can travel faster than the speed of light, while causality asserts that an event cannot precede its cause" " From biochemistry to physiology and pathology to pharmacology, each branch of medicine contributes its unique perspective, collaborating harmoniously to  paint a comprehensive picture of health and disease"

Replace-DocText " However, entanglement seems to violate these principles by allowing particles to instantaneously communicate their properties to their entangled partners, irrespective of the distance separating them. This perplexing phenomenon has fueled debates, challenging our notions of space, time, and the underlying fabric of reality" " As students embarking on this wondrous and demanding journey, we stand at the threshold of discovery, ready to delve into the intricacies of medicine, unraveling the mysteries that lie within"

Replace-DocText "The implications of quantum entanglement extend beyond the realm of theoretical physics" "In our pursuit of knowledge, medicine presents us with boundless opportunities for exploration"

Replace-DocText " Its potential applications span a wide range of fields, including cryptography, computing, and communication" " We will unravel the enigmas of genetic coding,Jie Pou Xue De Jing Yi , the intricate workings of the immune system"

Replace-DocText " harnessing the power of entangled particles could lead to the development of unbreakable codes, exponentially faster computers, and more secure communication networks" " We will witness the birth of medical innovations that transcend the boundaries of what we thought was possible"

Replace-DocText " As we continue to unravel the intricacies of this remarkable phenomenon, we stand at the precipice of a new era, poised to revolutionize our understanding of the universe and unlock unfathomable technological advancements" " From ground-breaking treatments for previously incurable diseases to the advent of preventative therapies, we will witness medicine's transformative impact on human lives"

# ---- Summary body paragraph ----
Replace-DocText "Quantum entanglement, a profound phenomenon in quantum physics, unveils the interconnectedness of particles beyond the constraints of space and time" "Through our study of medicine, we will gain an unwavering respect for the resilience of the human body, the complexity of disease processes, and the eternal pursuit of well-being"

Replace-DocText " Defying conventional notions of locality and causality, entanglement allows particles to instantaneously communicate their properties, regardless of their physical separation" " We will unravel the intricate pathways of disease, investigate the potent mechanisms of treatment, and ponder the profound nature of patient care"

Replace-DocText " This enigmatic phenomenon has ignited a revolution in our understanding of the universe, challenging fundamental concepts and opening up new frontiers of scientific exploration. With its potential applications spanning cryptography, computing, and communication, quantum entanglement holds the promise of transformative technologies that could revolutionize our world. As we delve deeper into its mysteries, we stand on the threshold of a new era, ready to unlock the extraordinary capabilities of this interconnected universe" " With each discovery, we draw inspiration from the enduring legacy of medical pioneers, walking in their footsteps as we continue the journey toward alleviating suffering and extending the boundaries of human life"

# ---- Append a new empty paragraph at the very end of the document (before sectPr) ----
$lastParaCount = $d.Paragraphs.Count
$endPos = $d.Paragraphs($lastParaCount).Range.End
$endRng = $d.Range($endPos, $endPos)
$endRng.Text = "`r"

Write-Host "Edits applied."
